# Generate Report for Archive
#
# 1. Update the "Status" value shown on every sheet from "Ready for handoff"
#    to "In Translation" (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the "Status" column on every sheet (Overview columns E & F;
#    zh-cn and de-de column C) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

# --- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
